$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Consumo")

# "Subarea -> Oficina 001": tidy up the "regla" text for CM-DS-001 / CM-DS-002
# (drop the redundant subarea number, keep the generic "x m2 = N mL" wording)
$ws.Range("D3").Value = "x m2 = 1 mL"
$ws.Range("D4").Value = "x m2 = 2 mL"

# Move the active selection, matching the saved view state in the workbook.
$ws.Range("E15").Select()
